# Update Flt3l-Flt3 LR-pairs data per Dr Hou's advice (re-run with updated permutation counts)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("ECs", "Flt3l", "Flt3", "M1", 3, 1, 5.918195, 17.754585, 0.2395278098354859, 0.2548411403336102, 3, 1, 2.890020333333334, 8.670061, 0.6127585034662664, 0.6135777502279139, 17.10370388663167, 153.933334979685, 0.1467727022933448, 0.1563648535514126)
for ($i = 0; $i -lt $row2.Length; $i++) { $ws.Cells.Item(2, $i + 1).Value = $row2[$i] }

$row3 = @("ECs", "Flt3l", "Flt3", "M2", 3, 1, 5.918195, 17.754585, 0.2395278098354859, 0.2548411403336102, 3, 1, 1.807497666666667, 5.422493, 0.3832359075370179, 0.3837482868421124, 10.69712365337833, 96.274112880405, 0.09179565758265669, 0.09779485101991324)
for ($i = 0; $i -lt $row3.Length; $i++) { $ws.Cells.Item(3, $i + 1).Value = $row3[$i] }

$row4 = @("ECs", "Flt3l", "Flt3", "sCs", 3, 1, 5.918195, 17.754585, 0.2395278098354859, 0.2548411403336102, 1, 0.5, 0.018892, 0.037784, 0.004005588996715721, 0.002673962929973791, 0.11180653994, 0.6708392396399999, 0.0009594499594844382, 0.0006814357622843222)
for ($i = 0; $i -lt $row4.Length; $i++) { $ws.Cells.Item(4, $i + 1).Value = $row4[$i] }

$row5 = @("FAPs", "Flt3l", "Flt3", "M1", 3, 1, 6.676481666666667, 20.029445, 0.2702180362464301, 0.287493433614434, 3, 1, 2.890020333333334, 8.670061, 0.6127585034662664, 0.6135777502279139, 19.29516777179389, 173.656509946145, 0.1655783994999559, 0.1763995742024426)
for ($i = 0; $i -lt $row5.Length; $i++) { $ws.Cells.Item(5, $i + 1).Value = $row5[$i] }

$row6 = @("FAPs", "Flt3l", "Flt3", "M2", 3, 1, 6.676481666666667, 20.029445, 0.2702180362464301, 0.287493433614434, 3, 1, 1.807497666666667, 5.422493, 0.3832359075370179, 0.3837482868421124, 12.06772503404278, 108.609525306385, 0.1035572543537714, 0.1103251126278956)
for ($i = 0; $i -lt $row6.Length; $i++) { $ws.Cells.Item(6, $i + 1).Value = $row6[$i] }

$row7 = @("FAPs", "Flt3l", "Flt3", "sCs", 3, 1, 6.676481666666667, 20.029445, 0.2702180362464301, 0.287493433614434, 1, 0.5, 0.018892, 0.037784, 0.004005588996715721, 0.002673962929973791, 0.1261320916466667, 0.7567925498799999, 0.001082382392702831, 0.0007687467840958776)
for ($i = 0; $i -lt $row7.Length; $i++) { $ws.Cells.Item(7, $i + 1).Value = $row7[$i] }

$row8 = @("M1", "Flt3l", "Flt3", "M1", 3, 1, 4.101432, 12.304296, 0.1659977449457439, 0.1766102008941509, 3, 1, 2.890020333333334, 8.670061, 0.6127585034662664, 0.6135777502279139, 11.853221875784, 106.678996882056, 0.101716529771729, 0.108364089731933)
for ($i = 0; $i -lt $row8.Length; $i++) { $ws.Cells.Item(8, $i + 1).Value = $row8[$i] }

$row9 = @("M1", "Flt3l", "Flt3", "M2", 3, 1, 4.101432, 12.304296, 0.1659977449457439, 0.1766102008941509, 3, 1, 1.807497666666667, 5.422493, 0.3832359075370179, 0.3837482868421124, 7.413328769992, 66.71995892992801, 0.06361629643338058, 0.06777386203197172)
for ($i = 0; $i -lt $row9.Length; $i++) { $ws.Cells.Item(9, $i + 1).Value = $row9[$i] }

$row10 = @("M1", "Flt3l", "Flt3", "sCs", 3, 1, 4.101432, 12.304296, 0.1659977449457439, 0.1766102008941509, 1, 0.5, 0.018892, 0.037784, 0.004005588996715721, 0.002673962929973791, 0.077484253344, 0.464905520064, 0.0006649187406342945, 0.0004722491302461836)
for ($i = 0; $i -lt $row10.Length; $i++) { $ws.Cells.Item(10, $i + 1).Value = $row10[$i] }

$row11 = @("M2", "Flt3l", "Flt3", "M1", 3, 1, 3.557602666666666, 10.672808, 0.143987275683135, 0.1531925731455665, 3, 1, 2.890020333333334, 8.670061, 0.6127585034662664, 0.6135777502279139, 10.28154404458756, 92.533896401288, 0.0882294275657825, 0.09399555438228183)
for ($i = 0; $i -lt $row11.Length; $i++) { $ws.Cells.Item(11, $i + 1).Value = $row11[$i] }

$row12 = @("M2", "Flt3l", "Flt3", "M2", 3, 1, 3.557602666666666, 10.672808, 0.143987275683135, 0.1531925731455665, 3, 1, 1.807497666666667, 5.422493, 0.3832359075370179, 0.3837482868421124, 6.430358518927111, 57.873226670344, 0.05518109427020901, 0.05878738750154613)
for ($i = 0; $i -lt $row12.Length; $i++) { $ws.Cells.Item(12, $i + 1).Value = $row12[$i] }

$row13 = @("M2", "Flt3l", "Flt3", "sCs", 3, 1, 3.557602666666666, 10.672808, 0.143987275683135, 0.1531925731455665, 1, 0.5, 0.018892, 0.037784, 0.004005588996715721, 0.002673962929973791, 0.06721022957866667, 0.403261377472, 0.0005767538471434386, 0.0004096312617385432)
for ($i = 0; $i -lt $row13.Length; $i++) { $ws.Cells.Item(13, $i + 1).Value = $row13[$i] }

$row14 = @("sCs", "Flt3l", "Flt3", "M1", 2, 1, 4.454046, 8.908092, 0.180269133289205, 0.1278626520122385, 3, 1, 2.890020333333334, 8.670061, 0.6127585034662664, 0.6135777502279139, 12.872283505602, 77.23370103361201, 0.1104614443354542, 0.07845367835984397)
for ($i = 0; $i -lt $row14.Length; $i++) { $ws.Cells.Item(14, $i + 1).Value = $row14[$i] }

$row15 = @("sCs", "Flt3l", "Flt3", "M2", 2, 1, 4.454046, 8.908092, 0.180269133289205, 0.1278626520122385, 3, 1, 1.807497666666667, 5.422493, 0.3832359075370179, 0.3837482868421124, 8.050677752226, 48.304066513356, 0.06908560489700014, 0.04906707366078571)
for ($i = 0; $i -lt $row15.Length; $i++) { $ws.Cells.Item(15, $i + 1).Value = $row15[$i] }

$row16 = @("sCs", "Flt3l", "Flt3", "sCs", 2, 1, 4.454046, 8.908092, 0.180269133289205, 0.1278626520122385, 1, 0.5, 0.018892, 0.037784, 0.004005588996715721, 0.002673962929973791, 0.084145837032, 0.336583348128, 0.0007220840567507195, 0.0003418999916088646)
for ($i = 0; $i -lt $row16.Length; $i++) { $ws.Cells.Item(16, $i + 1).Value = $row16[$i] }
